# Fruta / hortaliza, semanal
# Insert two new price rows (179 and 180) into the Camote/Zapallo data block,
# pushing the existing rows 179-208 down to 181-210.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 179; this shifts the old
# rows 179..208 down to 181..210 and keeps their formatting/styles intact.
$ws.Rows("179:180").Insert()

# --- New row 179 ---
$ws.Range("A179").Value = 7
$ws.Range("B179").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C179").Value = "Ñuble"
$ws.Range("D179").Value = 44889
$ws.Range("E179").Value = 16
$ws.Range("F179").Value = 100112045
$ws.Range("G179").Value = "Zapallo"
$ws.Range("H179").Value = "Camote"
$ws.Range("I179").Value = "1a nueva(o)"
$ws.Range("J179").Value = 300
$ws.Range("K179").Value = 1000
$ws.Range("L179").Value = 1100
$ws.Range("M179").Value = 1050
$ws.Range("N179").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O179").Value = "Perú"
$ws.Range("P179").Value = 1050
$ws.Range("Q179").Value = 1
$ws.Range("R179").Value = "Hortaliza"

# --- New row 180 ---
$ws.Range("A180").Value = 7
$ws.Range("B180").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C180").Value = "Ñuble"
$ws.Range("D180").Value = 44889
$ws.Range("E180").Value = 16
$ws.Range("F180").Value = 100112045
$ws.Range("G180").Value = "Zapallo"
$ws.Range("H180").Value = "Paine"
$ws.Range("I180").Value = "1a (guarda)"
$ws.Range("J180").Value = 600
$ws.Range("K180").Value = 550
$ws.Range("L180").Value = 600
$ws.Range("M180").Value = 575
$ws.Range("N180").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O180").Value = "Región del Maule"
$ws.Range("P180").Value = 575
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = "Hortaliza"
